$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 200.66667
$ws.Range("I4").Value = 200.66667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200.66667
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -86.66667000000001

$ws.Range("H40").Value = 6662.722
$ws.Range("I40").Value = 4825.8335
$ws.Range("J40").Value = 7581.1665
$ws.Range("K40").Value = 4825.8335
$ws.Range("L40").Value = 7581.1665
$ws.Range("M40").Value = -4650.8335
$ws.Range("N40").Value = -7931.1665

$ws.Range("H41").Value = 445.5
$ws.Range("I41").Value = 331.875
$ws.Range("J41").Value = 900
$ws.Range("K41").Value = 331.875
$ws.Range("L41").Value = 900
$ws.Range("M41").Value = 108.125
$ws.Range("N41").Value = -1780

$ws.Range("H70").Value = 3717.7646
$ws.Range("I70").Value = 2716.6667
$ws.Range("K70").Value = 8150.000100000001
$ws.Range("M70").Value = -7880.000100000001

$ws.Range("H73").Value = 3717.7646
$ws.Range("I73").Value = 2716.6667
$ws.Range("K73").Value = 8150.000100000001
$ws.Range("M73").Value = -7214.000100000001

$ws.Range("H98").Value = 2699.5
$ws.Range("I98").Value = 2699.5
$ws.Range("K98").Value = 2699.5
$ws.Range("M98").Value = -1201.5

$ws.Range("H116").Value = 4542.6665
$ws.Range("I116").Value = 3325
$ws.Range("J116").Value = 5516.8
$ws.Range("K116").Value = 3325
$ws.Range("L116").Value = 5516.8
$ws.Range("M116").Value = 117
$ws.Range("N116").Value = -12400.8

$ws.Range("H122").Value = 2699.5
$ws.Range("I122").Value = 2699.5
$ws.Range("K122").Value = 8098.5
$ws.Range("M122").Value = -5648.5

$ws.Range("H132").Value = 1397.3077
$ws.Range("I132").Value = 1257.4
$ws.Range("J132").Value = 1863.6666
$ws.Range("K132").Value = 3772.2
$ws.Range("L132").Value = 5590.9998
$ws.Range("M132").Value = -1242.2
$ws.Range("N132").Value = -10650.9998

$ws.Range("H137").Value = 1966.6666
$ws.Range("I137").Value = 1950
$ws.Range("K137").Value = 5850
$ws.Range("M137").Value = -3300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4003.0293
$ws.Range("I2").Value = 3671.7144
$ws.Range("K2").Value = 3671.7144
$ws.Range("M2").Value = -3558.7144

$ws.Range("H32").Value = 5688.4565
$ws.Range("I32").Value = 4964.6587
$ws.Range("J32").Value = 11623.6
$ws.Range("K32").Value = 4964.6587
$ws.Range("L32").Value = 11623.6
$ws.Range("M32").Value = -4677.6587
$ws.Range("N32").Value = -12197.6

$ws.Range("H45").Value = 9500
$ws.Range("I45").Value = 7500
$ws.Range("K45").Value = 7500
$ws.Range("M45").Value = -7123

$ws.Range("H116").Value = 4003.0293
$ws.Range("I116").Value = 3671.7144
$ws.Range("K116").Value = 3671.7144
$ws.Range("M116").Value = -1377.7144

$ws.Range("H132").Value = 3787.4146
$ws.Range("I132").Value = 3196.8333
$ws.Range("K132").Value = 9590.499899999999
$ws.Range("M132").Value = -7060.499899999999

$ws.Range("H139").Value = 97809.664
$ws.Range("J139").Value = 97809.664
$ws.Range("L139").Value = 97809.664
$ws.Range("N139").Value = -108089.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4003.0293
$ws.Range("I3").Value = 3671.7144
$ws.Range("K3").Value = 3671.7144
$ws.Range("M3").Value = -3557.7144

$ws.Range("H86").Value = 18186884
$ws.Range("I86").Value = 6619.5
$ws.Range("K86").Value = 6619.5
$ws.Range("M86").Value = -5496.5

$ws.Range("H89").Value = 18186884
$ws.Range("I89").Value = 6619.5
$ws.Range("K89").Value = 33097.5
$ws.Range("M89").Value = -27481.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3461.9592
$ws.Range("I31").Value = 2647.4138
$ws.Range("K31").Value = 2647.4138
$ws.Range("M31").Value = -2352.4138

$ws.Range("H34").Value = 3461.9592
$ws.Range("I34").Value = 2647.4138
$ws.Range("K34").Value = 2647.4138
$ws.Range("M34").Value = -2445.4138

$ws.Range("H58").Value = 4547.45
$ws.Range("I58").Value = 1926
$ws.Range("K58").Value = 1926
$ws.Range("M58").Value = -1723

$ws.Range("H122").Value = 4131.6875
$ws.Range("I122").Value = 4042.3572
$ws.Range("K122").Value = 12127.0716
$ws.Range("M122").Value = -9677.071599999999

$ws.Range("H134").Value = 6127.7666
$ws.Range("I134").Value = 2544.6316
$ws.Range("K134").Value = 7633.8948
$ws.Range("M134").Value = -5098.8948

$ws.Range("H136").Value = 4547.45
$ws.Range("I136").Value = 1926
$ws.Range("K136").Value = 5778
$ws.Range("M136").Value = -3228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2093.4827
$ws.Range("I5").Value = 1899.4
$ws.Range("J5").Value = 2301.4285
$ws.Range("K5").Value = 5698.200000000001
$ws.Range("L5").Value = 6904.2855
$ws.Range("M5").Value = -5586.200000000001
$ws.Range("N5").Value = -7128.2855

$ws.Range("H34").Value = 725.2
$ws.Range("I34").Value = 264.14285
$ws.Range("J34").Value = 1801
$ws.Range("K34").Value = 792.4285500000001
$ws.Range("L34").Value = 5403
$ws.Range("M34").Value = -708.4285500000001
$ws.Range("N34").Value = -5571

$ws.Range("H68").Value = 31251946
$ws.Range("I68").Value = 38463508
$ws.Range("J68").Value = 1847.3334
$ws.Range("K68").Value = 115390524
$ws.Range("L68").Value = 5542.0002
$ws.Range("M68").Value = -115389713
$ws.Range("N68").Value = -7164.0002

$ws.Range("H71").Value = 31251946
$ws.Range("I71").Value = 38463508
$ws.Range("J71").Value = 1847.3334
$ws.Range("K71").Value = 346171572
$ws.Range("L71").Value = 16626.0006
$ws.Range("M71").Value = -346167516
$ws.Range("N71").Value = -24738.0006

$ws.Range("H135").Value = 2093.4827
$ws.Range("I135").Value = 1899.4
$ws.Range("J135").Value = 2301.4285
$ws.Range("K135").Value = 17094.6
$ws.Range("L135").Value = 20712.8565
$ws.Range("M135").Value = -14559.6
$ws.Range("N135").Value = -25782.8565

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3313.5
$ws.Range("I80").Value = 3805
$ws.Range("J80").Value = 3275.6924
$ws.Range("K80").Value = 3805
$ws.Range("L80").Value = 3275.6924
$ws.Range("M80").Value = -2807
$ws.Range("N80").Value = -5271.6924

$ws.Range("H83").Value = 3313.5
$ws.Range("I83").Value = 3805
$ws.Range("J83").Value = 3275.6924
$ws.Range("K83").Value = 19025
$ws.Range("L83").Value = 16378.462
$ws.Range("M83").Value = -14033
$ws.Range("N83").Value = -26362.462

$ws.Range("H102").Value = 2627.2
$ws.Range("I102").Value = 2109.1
$ws.Range("K102").Value = 2109.1
$ws.Range("M102").Value = -487.0999999999999

$ws.Range("H132").Value = 1600
$ws.Range("I132").Value = 1550
$ws.Range("K132").Value = 4650
$ws.Range("M132").Value = -2120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2815.1667
$ws.Range("I16").Value = 2815.1667
$ws.Range("K16").Value = 2815.1667
$ws.Range("M16").Value = -2645.1667

$ws.Range("H40").Value = 3554.182
$ws.Range("I40").Value = 2646.1333
$ws.Range("K40").Value = 2646.1333
$ws.Range("M40").Value = -2510.1333

$ws.Range("H44").Value = 22955
$ws.Range("I44").Value = 22954
$ws.Range("J44").Value = 22956
$ws.Range("K44").Value = 22954
$ws.Range("L44").Value = 22956
$ws.Range("M44").Value = -22498
$ws.Range("N44").Value = -23868

$ws.Range("H122").Value = 4848.5
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13083.714
$ws.Range("J41").Value = 14433.167
$ws.Range("L41").Value = 14433.167
$ws.Range("N41").Value = -15213.167

$ws.Range("H70").Value = 33666.668
$ws.Range("I70").Value = 20000
$ws.Range("K70").Value = 20000
$ws.Range("M70").Value = -19685

$ws.Range("H73").Value = 33666.668
$ws.Range("I73").Value = 20000
$ws.Range("K73").Value = 20000
$ws.Range("M73").Value = -18908

$ws.Range("H136").Value = 4345.488
$ws.Range("I136").Value = 3192.8438
$ws.Range("J136").Value = 8443.777
$ws.Range("K136").Value = 9578.5314
$ws.Range("L136").Value = 25331.331
$ws.Range("M136").Value = -7028.5314
$ws.Range("N136").Value = -30431.331
